# Apply "Update countries & provincias Spain" data refresh to the Pais sheet.
# Source data refreshed as of 5 de Agosto de 2020 a las 02:57 (was 01:40).
# A handful of countries swapped rank (column A) because their updated
# "Casos totales" (column B) overtook their former neighbour's.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp row
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 02:57"

# Estados Unidos (row 4) - updated totals, no rank change
$ws.Range("B4").Value = 4918406
$ws.Range("C4").Value = 54490
$ws.Range("D4").Value = 2481223
$ws.Range("E4").Value = 2276893
$ws.Range("G4").Value = 1362
$ws.Range("H4").Value = 160290

# Peru (row 10) - updated totals, no rank change
$ws.Range("B10").Value = 439890
$ws.Range("C10").Value = 6790
$ws.Range("D10").Value = 302457
$ws.Range("E10").Value = 117426
$ws.Range("G10").Value = 196
$ws.Range("H10").Value = 20007

# Canada (row 25) - updated totals, no rank change
$ws.Range("D25").Value = 102450
$ws.Range("E25").Value = 6384

# Kuwait / Panama swap rank (rows 40-41): Panama's refreshed numbers
# overtake Kuwait, so Panama moves up to row 40 and Kuwait drops to row 41.
$ws.Range("A40").Value = "Panama"
$ws.Range("B40").Value = 69424
$ws.Range("C40").Value = 968
$ws.Range("D40").Value = 43330
$ws.Range("E40").Value = 24572
$ws.Range("G40").Value = 25
$ws.Range("H40").Value = 1522

$ws.Range("A41").Value = "Kuwait"
$ws.Range("B41").Value = 68774
$ws.Range("C41").Value = 475
$ws.Range("D41").Value = 60326
$ws.Range("E41").Value = 7983
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 465

# Chequia / Camerun swap rank (rows 74-75): Camerun moves up to row 74.
$ws.Range("A74").Value = "Camerun"
$ws.Range("B74").Value = 17718
$ws.Range("C74").Value = 463
$ws.Range("D74").Value = 15320
$ws.Range("E74").Value = 2007
$ws.Range("H74").Value = 391

$ws.Range("A75").Value = "Chequia"
$ws.Range("B75").Value = 17286
$ws.Range("C75").Value = 278
$ws.Range("D75").Value = 11812
$ws.Range("E75").Value = 5091
$ws.Range("H75").Value = 383

# Albania / Paraguay swap rank (rows 98-99): Paraguay moves up to row 98.
$ws.Range("A98").Value = "Paraguay"
$ws.Range("B98").Value = 5852
$ws.Range("C98").Value = 128
$ws.Range("D98").Value = 4645
$ws.Range("E98").Value = 1148
$ws.Range("H98").Value = 59

$ws.Range("A99").Value = "Albania"
$ws.Range("B99").Value = 5750
$ws.Range("C99").Value = 130
$ws.Range("D99").Value = 3031
$ws.Range("E99").Value = 2543
$ws.Range("G99").Value = 4
$ws.Range("H99").Value = 176

# Zimbabue / Libia swap rank (rows 109-110): Libia moves up to row 109.
$ws.Range("A109").Value = "Libia"
$ws.Range("B109").Value = 4224
$ws.Range("C109").Value = 161
$ws.Range("D109").Value = 633
$ws.Range("E109").Value = 3495
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = 96

$ws.Range("A110").Value = "Zimbabue"
$ws.Range("B110").Value = 4221
$ws.Range("C110").Value = 146
$ws.Range("D110").Value = 1238
$ws.Range("E110").Value = 2902
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 81

# Burkina Faso (row 147) - updated totals, no rank change
$ws.Range("E147").Value = 152
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 54

# Reunion (row 163) - updated totals, no rank change
$ws.Range("B163").Value = 669
$ws.Range("C163").Value = 2
$ws.Range("E163").Value = 73

# Trinidad yTobago (row 179) - updated totals, no rank change
$ws.Range("B179").Value = 194
$ws.Range("C179").Value = 12
$ws.Range("E179").Value = 51
